$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'63.647.34"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.33%  '

$c = $ws.Range("D3")
$c.Value = "'2.655.33"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.90%  '

$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$c = $ws.Range("D5")
$c.Value = "'591.73"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.51%  '

$c = $ws.Range("D6")
$c.Value = "'144.85"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.84%  '

$c = $ws.Range("D7")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '

$c = $ws.Range("D8")
$c.Value = "'0.589"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.73%  '

$c = $ws.Range("D9")
$c.Value = "'2.654.15"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.88%  '

$ws.Range("E10").Value = '  -0.26%  '

$c = $ws.Range("D11")
$c.Value = "'5.62"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("E12").Value = '  +0.49%  '

$c = $ws.Range("D13")
$c.Value = "'0.354"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.50%  '

$c = $ws.Range("D14")
$c.Value = "'27.53"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '

$c = $ws.Range("D15")
$c.Value = "'3.124.89"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.63%  '

$c = $ws.Range("D16")
$c.Value = "'63.517.80"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.25%  '

$c = $ws.Range("D17")
$c.Value = "'0.0000146"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.09%  '

$c = $ws.Range("D18")
$c.Value = "'2.631.97"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.86%  '

$c = $ws.Range("D19")
$c.Value = "'11.44"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.64%  '

$c = $ws.Range("D20")
$c.Value = "'341.00"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.11%  '

$c = $ws.Range("D21")
$c.Value = "'4.37"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.16%  '

$c = $ws.Range("D22")
$c.Value = "'6.73"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.02%  '

$ws.Range("E23").Value = '  +0.20%  '

$c = $ws.Range("D24")
$c.Value = "'67.65"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.47%  '

$c = $ws.Range("D25")
$c.Value = "'1.64"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.76%  '

$c = $ws.Range("D26")
$c.Value = "'1.56"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +7.07%  '

$ws.Range("E27").Value = '  -0.11%  '

$c = $ws.Range("D28")
$c.Value = "'553.31"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +18.66%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D29")
$c.Value = "'8.47"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.45%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D30")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '

$c = $ws.Range("D31")
$c.Value = "'7.79"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.82%  '

$c = $ws.Range("D32")
$c.Value = "'1.83"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +14.38%  '

$c = $ws.Range("D33")
$c.Value = "'1.98"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +3.14%  '

$c = $ws.Range("D34")
$c.Value = "'0.0₃0812"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.50%  '

$c = $ws.Range("D35")
$c.Value = "'174.24"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.27%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D36")
$c.Value = "'4.88"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +8.83%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D37")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '

$c = $ws.Range("D38")
$c.Value = "'0.403"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.14%  '

$c = $ws.Range("D39")
$c.Value = "'19.12"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.87%  '

$c = $ws.Range("D40")
$c.Value = "'1.82"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +7.56%  '

$c = $ws.Range("D41")
$c.Value = "'170.53"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +7.99%  '

$c = $ws.Range("D42")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.06%  '

$c = $ws.Range("D43")
$c.Value = "'40.33"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.09%  '

$c = $ws.Range("D44")
$c.Value = "'3.75"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.32%  '

$c = $ws.Range("D45")
$c.Value = "'22.36"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +6.82%  '

$c = $ws.Range("D46")
$c.Value = "'0.631"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.26%  '

$c = $ws.Range("D47")
$c.Value = "'0.0554"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.95%  '

$c = $ws.Range("D48")
$c.Value = "'0.0962"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.34%  '

$c = $ws.Range("D49")
$c.Value = "'0.0239"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.59%  '

$c = $ws.Range("D50")
$c.Value = "'18.76"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.77%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D51")
$c.Value = "'1.71"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.15%  '
